$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Password value in B2 (was "ngasce@775", now "Ngasce@123")
$ws.Range("B2").Value = "Ngasce@123"

# Remove the extra row of student/login details (row 3: 54443434 / grhru)
$ws.Rows("3").Delete()

# Move the active selection to B2
$ws.Range("B2").Select()
